# Implement account bills rules handler
# Updates the "Title" column (G) entries to include the matched rule's
# output (account/category name appended to the original title), widens
# column G to fit the new, longer text, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Title values for each transaction row with the rule output.
$ws.Range("G2").Value = "Title: 0090729 in ABC"
$ws.Range("G3").Value = "Title: 0090744 in ABC"
$ws.Range("G4").Value = "Title: 00091840 in Little Shop"

# Widen column G (Title) so the longer values are fully visible.
$ws.Columns.Item(7).ColumnWidth = 26

# Move / restore the active selection.
$ws.Range("G5").Select()
